$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PA73405 - Attrition by Job 2009")

# Row 7: Department Leader
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "2009"
$ws.Range("A7").Style = "Normal"
$ws.Range("B7").Value = "Department Leader"
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 3
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "0.0%"
$ws.Range("F7").Style = "Normal"

# Row 8: Regional Leader
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "2009"
$ws.Range("A8").Style = "Normal"
$ws.Range("B8").Value = "Regional Leader"
$ws.Range("C8").Value = 9
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 9
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "0.0%"
$ws.Range("F8").Style = "Normal"

# Resize the Table3 ListObject to include the new rows
$table = $ws.ListObjects.Item("Table3")
$table.Resize($ws.Range("A1:F8"))
